$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (H) values updated
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 478
$wsOff.Range("C2").Value = 341
$wsOff.Range("D2").Value = 133
$wsOff.Range("E2").Value = 66
$wsOff.Range("F2").Value = 8

# Sheet "DEF" - row 2 (H) values updated
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 544
$wsDef.Range("C2").Value = 363
$wsDef.Range("D2").Value = 127
$wsDef.Range("E2").Value = 48
$wsDef.Range("F2").Value = 10
